$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 14762
$ws.Range("F5").Value = 1669
$ws.Range("F6").Value = 510
$ws.Range("F7").Value = 2138
$ws.Range("F8").Value = 1312
$ws.Range("F9").Value = 1994
$ws.Range("F11").Value = 45
$ws.Range("F12").Value = 2363
$ws.Range("F13").Value = 625
$ws.Range("F15").Value = 3687
$ws.Range("F17").Value = 349
$ws.Range("F18").Value = 2740
$ws.Range("F19").Value = 708
$ws.Range("F22").Value = 1944
$ws.Range("F23").Value = 1142
$ws.Range("F24").Value = 1667
$ws.Range("F25").Value = 348
$ws.Range("F26").Value = 181
$ws.Range("F27").Value = 7688
$ws.Range("F28").Value = 5305
$ws.Range("F29").Value = 336
$ws.Range("F31").Value = 730
$ws.Range("F32").Value = 737
$ws.Range("F33").Value = 3422
$ws.Range("F35").Value = 933
$ws.Range("F36").Value = 366
$ws.Range("F37").Value = 160
$ws.Range("F38").Value = 130
$ws.Range("F39").Value = 4527
$ws.Range("F40").Value = 755
$ws.Range("F41").Value = 39
$ws.Range("F42").Value = 356

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F18").Value = 126

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 8079
$ws.Range("F3").Value = 327
$ws.Range("F4").Value = 1163

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 8079
$ws.Range("F4").Value = 327
$ws.Range("F5").Value = 1163
$ws.Range("F6").Value = 14762
$ws.Range("F9").Value = 1669
$ws.Range("F10").Value = 510
$ws.Range("F11").Value = 1312
$ws.Range("F12").Value = 1994
$ws.Range("F14").Value = 45
$ws.Range("F15").Value = 625
$ws.Range("F17").Value = 3687
$ws.Range("F18").Value = 349
$ws.Range("F19").Value = 2740
$ws.Range("F20").Value = 708
$ws.Range("F23").Value = 1944
$ws.Range("F29").Value = 1667
$ws.Range("F31").Value = 348
$ws.Range("F32").Value = 181
$ws.Range("F33").Value = 7688
$ws.Range("F34").Value = 5305
$ws.Range("F35").Value = 336
$ws.Range("F36").Value = 730
$ws.Range("F37").Value = 3422
$ws.Range("F39").Value = 933
$ws.Range("F40").Value = 366
$ws.Range("F42").Value = 130
$ws.Range("F43").Value = 4527
$ws.Range("F44").Value = 755
$ws.Range("F45").Value = 39
$ws.Range("F46").Value = 356
